# Tidsregistrering i PTE projektet Simon Nielsen.xlsx
# Apply: "Har lavet SSD for UC6 / Tidsregistrering opdateret"
# Fills in the timesheet rows for 8/3, 10/3, 13/3, 14/3 and 15/3 2017,
# widens column F (and sizes the two new helper columns G/H), and moves
# the selection/cursor down to where the user left off (F33).
#
# NB: this interpreter only binds function parameters positionally, so
# Set-Entry takes its arguments in a fixed order and uses $null for any
# field that row doesn't have.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

function Set-Entry($Row, $Date, $Role, $Activity, $Start, $End, $Duration) {
    if ($null -ne $Date) {
        $c = $ws.Cells.Item($Row, 1)
        $c.Value = $Date
        $c.NumberFormat = "m/d/yy"
    }
    if ($null -ne $Role) {
        $ws.Cells.Item($Row, 5).Value = $Role
    }
    if ($null -ne $Activity) {
        $ws.Cells.Item($Row, 6).Value = $Activity
    }
    if ($null -ne $Start) {
        $c = $ws.Cells.Item($Row, 7)
        $c.Value = $Start
        $c.NumberFormat = "h:mm"
    }
    if ($null -ne $End) {
        $c = $ws.Cells.Item($Row, 8)
        $c.Value = $End
        $c.NumberFormat = "h:mm"
    }
    if ($null -ne $Duration) {
        $ws.Cells.Item($Row, 9).Value = $Duration
    }
}

# --- Onsdag 8-3-2017 ------------------------------------------------------
Set-Entry 13 42802 "Any Role" "Rettet & samlet Domain Model" 0.45833333333333331 0.52083333333333337 "1.5t"
Set-Entry 14 $null "System Analyst" "Operationskontrakt 002" 0.54166666666666663 0.60416666666666663 "1.5t"
Set-Entry 15 $null $null $null $null $null "3t"

# --- Fredag 10-3-2017 ------------------------------------------------------
Set-Entry 16 42804 "Any Role" "Fælles design OC3" 0.34375 0.42708333333333331 "2t"
Set-Entry 17 $null "Reviewer" "Review OC2" 0.4375 0.47916666666666669 "1t"
Set-Entry 18 $null $null $null $null $null "3t"

# --- Mandag 13-3-2017 / Tirsdag 14-3-2017 : fraværende ---------------------
Set-Entry 20 42807 $null "Fraværende" $null $null $null
Set-Entry 21 42808 $null "Fraværende" $null $null $null

# --- Onsdag 15-3-2017 -------------------------------------------------------
Set-Entry 23 42809 "Requirements Specifier" "Test Case OC6" 0.35069444444444442 0.40972222222222227 "1t 25m"
Set-Entry 24 $null "Requirements Specifier" "SSD UC3" 0.41666666666666669 0.4375 "30m"
Set-Entry 25 $null "Reviewer" "Review af OC6 design" 0.4548611111111111 0.47222222222222227 $null
Set-Entry 26 $null "Reviewer" "Review af SSD3" 0.47222222222222227 0.4770833333333333 $null
Set-Entry 27 $null "Implementer" "Implementeret test case for OC 6" 0.50347222222222221 0.54513888888888895 $null
Set-Entry 28 $null "Implementer" "Implementeret design for OC6" 0.54513888888888895 0.57222222222222219 $null
Set-Entry 29 $null $null "Lavet SSD for UC9" 0.57638888888888895 0.58333333333333337 $null
Set-Entry 30 $null "Reviewer" "Krydscheck for UC6" 0.59375 0.62847222222222221 $null
Set-Entry 31 $null "Requirements Specifier" "SSD UC6" 0.63541666666666663 0.64583333333333337 $null

# --- Column sizing: widen the activity column, and size the two new
#     helper columns that now sit between "Aktivitet" and "Starttid" ------
$ws.Columns.Item(6).ColumnWidth = 48.04
$ws.Columns.Item(7).ColumnWidth = 6.46
$ws.Columns.Item(8).ColumnWidth = 7.2

# --- Cursor / selection: leave it where the user stopped typing -----------
$ws.Activate()
$ws.Range("F33").Select()
